$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '69.587.33'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '3.887.77'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.83'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.97'
$ws.Range('E6').Value = '  +4.42%  '
$ws.Range('D7').Value = '3.889.22'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('E13').Value = '  +4.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.21'
$ws.Range('E14').Value = '  +3.46%  '
$ws.Range('D15').Value = '4.546.69'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').Value = '3.898.80'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').Value = '69.651.23'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.73'
$ws.Range('E18').Value = '  +9.12%  '
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.04'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '489.95'
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.745'
$ws.Range('E23').Value = '  +3.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000165'
$ws.Range('E24').Value = '  +3.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.34'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.33'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.13'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = '4.041.66'
$ws.Range('E31').Value = '  +1.18%  '
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.83'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.88'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('D35').Value = '3.858.69'
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.144'
$ws.Range('E37').Value = '  +2.35%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.41'
$ws.Range('E38').Value = '  +14.63%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.12'
$ws.Range('E39').Value = '  +3.84%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +2.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.08'
$ws.Range('E43').Value = '  +4.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '435.27'
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '48.08'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000276'
$ws.Range('E48').Value = '  +21.25%  '
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '40.19'
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '142.25'
$ws.Range('E51').Value = '  -0.31%  '
